$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose records were dropped from the missing-data
# sample: original "RM 232" (row 26) and, after that shift, the original
# "SC 92" row (now row 27). Everything below shifts up to close the gaps.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-point / re-roll the remaining missing-data (NA) cells and values for
# this seed so the sheet matches the new error-calculation pass.
# A leading apostrophe is Excel's text-prefix marker, so assigning "'" is
# how a blank/NA marker cell (empty text, not a truly blank cell) is set;
# the style is put back to Normal afterwards so only the value changes.
$ws.Range("C2").Value = 14.9
$ws.Range("E2").Value = -7.2
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = -15.4
$ws.Range("D5").Value = "'"
$ws.Range("D5").Style = "Normal"
$ws.Range("C6").Value = "'"
$ws.Range("C6").Style = "Normal"
$ws.Range("D8").Value = "'"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'"
$ws.Range("E8").Style = "Normal"
$ws.Range("C12").Value = 12.5
$ws.Range("E13").Value = -5.3
$ws.Range("C14").Value = "'"
$ws.Range("C14").Style = "Normal"
$ws.Range("D15").Value = -15.2
$ws.Range("D18").Value = -15.2
$ws.Range("D19").Value = "'"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'"
$ws.Range("E19").Style = "Normal"
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("D22").Value = "'"
$ws.Range("D22").Style = "Normal"
$ws.Range("C23").Value = "'"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = -13.9
$ws.Range("C24").Value = "'"
$ws.Range("C24").Style = "Normal"
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = "'"
$ws.Range("B27").Style = "Normal"
$ws.Range("D27").Value = "'"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = -5.9
$ws.Range("B30").Value = -19.7
$ws.Range("C31").Value = 15.3
$ws.Range("E31").Value = "'"
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'"
$ws.Range("B32").Style = "Normal"
$ws.Range("E32").Value = -6.4
$ws.Range("C33").Value = 10.4
